# Update the legacy GSC export "Chart" sheet with refreshed Video-Indexing data.
# The data window rolled forward by 3 days: the earliest 3 rows (2025-09-13,
# blank, 2025-09-14) dropped off, and 3 new days (2025-12-09..2025-12-11) were
# appended, shifting/growing the table from A1:D88 to A1:D89.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Chart")

# Column A date labels, 2025-09-15 through 2025-12-11 (88 rows, rows 2..89),
# spelled out explicitly so the values are not subject to any locale/clock
# dependent date parsing.
$dates = @("2025-09-15","2025-09-16","2025-09-17","2025-09-18","2025-09-19","2025-09-20","2025-09-21","2025-09-22","2025-09-23","2025-09-24","2025-09-25","2025-09-26","2025-09-27","2025-09-28","2025-09-29","2025-09-30","2025-10-01","2025-10-02","2025-10-03","2025-10-04","2025-10-05","2025-10-06","2025-10-07","2025-10-08","2025-10-09","2025-10-10","2025-10-11","2025-10-12","2025-10-13","2025-10-14","2025-10-15","2025-10-16","2025-10-17","2025-10-18","2025-10-19","2025-10-20","2025-10-21","2025-10-22","2025-10-23","2025-10-24","2025-10-25","2025-10-26","2025-10-27","2025-10-28","2025-10-29","2025-10-30","2025-10-31","2025-11-01","2025-11-02","2025-11-03","2025-11-04","2025-11-05","2025-11-06","2025-11-07","2025-11-08","2025-11-09","2025-11-10","2025-11-11","2025-11-12","2025-11-13","2025-11-14","2025-11-15","2025-11-16","2025-11-17","2025-11-18","2025-11-19","2025-11-20","2025-11-21","2025-11-22","2025-11-23","2025-11-24","2025-11-25","2025-11-26","2025-11-27","2025-11-28","2025-11-29","2025-11-30","2025-12-01","2025-12-02","2025-12-03","2025-12-04","2025-12-05","2025-12-06","2025-12-07","2025-12-08","2025-12-09","2025-12-10","2025-12-11")

# "Video indexed" (column B) step values for each day above. Columns C
# ("Impressions" per the header, always 0 in this export) and D are both 0
# for every row.
$bValues = @(11,11,11,11,16,16,16,16,17,17,17,18,18,18,18,19,19,19,19,19,19,19,19,19,19,19,19,19,19,19,19,19,19,19,19,19,17,17,17,14,14,14,14,14,14,14,13,13,13,13,11,11,11,9,9,9,9,9,9,9,8,8,8,8,8,8,8,6,6,6,6,6,6,6,6,6,6,6,6,6,6,4,4,4,4,4,4,4)

$lastRow = $bValues.Count + 1

# Force column A to text BEFORE writing the date strings so Excel stores
# them as literal text (matching the source export) instead of converting
# them to date serial numbers.
$dateRange = $ws.Range("A2:A$lastRow")
$dateRange.NumberFormat = "@"

for ($i = 0; $i -lt $bValues.Count; $i++) {
    $row = $i + 2

    $ws.Cells.Item($row, 1).Value = $dates[$i]
    $ws.Cells.Item($row, 2).Value = $bValues[$i]
    $ws.Cells.Item($row, 3).Value = 0
    $ws.Cells.Item($row, 4).Value = 0
}

# Drop the explicit "text" number-format we applied above so the cells fall
# back to the workbook's default (General) style, matching the original
# export's formatting.
$dateRange.ClearFormats()
